$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names in column A (rows 5-8 swap order/pairing) and point totals in
# column B (week 4 scores added in, formulas replaced with plain values,
# rows re-sorted descending by total points).

$ws.Range("A2").Value = "James Manley"
$ws.Range("B2").Value = 1384.8

$ws.Range("A3").Value = "Matt Piunti"
$ws.Range("B3").Value = 1266.3

$ws.Range("A4").Value = "Steven Carter"
$ws.Range("B4").Value = 1010.6

$ws.Range("A5").Value = "Todd Vinsant"
$ws.Range("B5").Value = 954.6

$ws.Range("A6").Value = "Jeremiah Gaddy"
$ws.Range("B6").Value = 845.6

$ws.Range("A7").Value = "Senay Semere"
$ws.Range("B7").Value = 788.6

$ws.Range("A8").Value = "Philip Milam"
$ws.Range("B8").Value = 671.2

$ws.Range("A9").Value = "Josh Lance"
$ws.Range("B9").Value = 671.2

$ws.Range("A10").Value = "Andrew Harrell"
$ws.Range("B10").Value = 250.8

$ws.Range("A11").Value = "Brandon Greife"
$ws.Range("B11").Value = 79

# Move the active selection as it appears in the saved file
$ws.Range("B14").Select()
